$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sims and main analysis")

# New output path/job on Eddie instead of the scratch dir (files >1 month old get purged there).
# Row 7: record the corrected re-analysis start date, job ID, completion date and comment.
$ws.Range("L7").Value = 45791
$ws.Range("L7").NumberFormat = "d-mmm-yy"
$ws.Range("L7").Interior.Color = 65535

$ws.Range("M7").Value = "48690588 (Eddie)"

$ws.Range("N7").Value = 45791
$ws.Range("N7").NumberFormat = "d-mmm-yy"
$ws.Range("N7").Interior.Color = 65535

$ws.Range("O7").Value = "NA"

# Row 13: record the corrected re-analysis start date.
$ws.Range("L13").Value = 45792
$ws.Range("L13").NumberFormat = "d-mmm-yy"
$ws.Range("L13").Interior.Color = 65535

# Move the cursor/selection off the old M23 spot, onto A13, and reset the scrolled view.
$ws.Activate()
$ws.Range("A13").Select()
